# issue #5: stock data from json to db
# Adds "category", "source_file" and "index" columns to the 股票 (stock) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column I ("category") - shifts old I(date)/J(legislator_name)/K(legislator_id) one to the right.
$ws.Columns.Item(9).Insert()

# Insert two new columns at M:N ("source_file", "index") after the (now shifted) legislator_id column L.
$ws.Range("M1:N1").EntireColumn.Insert()

# Header row
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Data rows: r = 2..13, A column already holds the original row index (72..83)
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp22e71"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value()
}
